# Applies the "Updated cryptos list" data refresh described in the commit.
# Price (col D) and Volume(1h) (col E) text values are updated for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value to a cell without letting Excel
# auto-convert numeric-looking text (e.g. "580.64") into a real number,
# then restore the cell to its original (unstyled) "Normal" style so the
# formatting/style table is left untouched.
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "66.831.26"
$ws.Cells.Item(2, 5).Value = "  -2.66%  "
$ws.Cells.Item(3, 4).Value = "2.457.75"
$ws.Cells.Item(3, 5).Value = "  -3.80%  "
$ws.Cells.Item(4, 5).Value = "  +0.01%  "
Set-TextCell 5 4 "580.64"
$ws.Cells.Item(5, 5).Value = "  -2.39%  "
Set-TextCell 6 4 "168.07"
$ws.Cells.Item(6, 5).Value = "  -4.92%  "
$ws.Cells.Item(7, 5).Value = "  +0.08%  "
Set-TextCell 8 4 "0.510"
$ws.Cells.Item(8, 5).Value = "  -3.19%  "
$ws.Cells.Item(9, 4).Value = "2.459.97"
$ws.Cells.Item(9, 5).Value = "  -3.69%  "
Set-TextCell 10 4 "0.133"
$ws.Cells.Item(10, 5).Value = "  -4.27%  "
$ws.Cells.Item(11, 5).Value = "  -0.99%  "
Set-TextCell 12 4 "4.88"
$ws.Cells.Item(12, 5).Value = "  -3.13%  "
Set-TextCell 13 4 "0.326"
$ws.Cells.Item(13, 5).Value = "  -5.64%  "
$ws.Cells.Item(14, 4).Value = "2.907.29"
$ws.Cells.Item(14, 5).Value = "  -1.68%  "
Set-TextCell 15 4 "25.24"
$ws.Cells.Item(15, 5).Value = "  -5.32%  "
$ws.Cells.Item(16, 4).Value = "66.744.80"
$ws.Cells.Item(16, 5).Value = "  -2.73%  "
Set-TextCell 17 4 "0.0000167"
$ws.Cells.Item(17, 5).Value = "  -5.99%  "
$ws.Cells.Item(18, 4).Value = "2.459.66"
$ws.Cells.Item(18, 5).Value = "  -3.54%  "
Set-TextCell 19 4 "10.94"
$ws.Cells.Item(19, 5).Value = "  -8.73%  "
Set-TextCell 20 4 "7.41"
$ws.Cells.Item(20, 5).Value = "  -8.11%  "
Set-TextCell 21 4 "349.30"
$ws.Cells.Item(21, 5).Value = "  -6.09%  "
Set-TextCell 22 4 "3.99"
$ws.Cells.Item(22, 5).Value = "  -4.45%  "
$ws.Cells.Item(23, 5).Value = "  -0.05%  "
Set-TextCell 24 4 "68.61"
$ws.Cells.Item(24, 5).Value = "  -4.98%  "
Set-TextCell 25 4 "4.19"
$ws.Cells.Item(25, 5).Value = "  -8.84%  "
Set-TextCell 26 4 "1.81"
$ws.Cells.Item(26, 5).Value = "  -5.74%  "
Set-TextCell 27 4 "9.08"
$ws.Cells.Item(27, 5).Value = "  -9.06%  "
Set-TextCell 28 4 "0.997"
$ws.Cells.Item(28, 5).Value = "  -49.47%  "
$ws.Cells.Item(29, 4).Value = "2.596.17"
$ws.Cells.Item(29, 5).Value = "  -2.77%  "
$ws.Cells.Item(30, 4).Value = "0.0₃0891"
$ws.Cells.Item(30, 5).Value = "  -8.73%  "
Set-TextCell 31 4 "505.64"
$ws.Cells.Item(31, 5).Value = "  -6.31%  "
Set-TextCell 32 4 "7.60"
$ws.Cells.Item(32, 5).Value = "  -8.83%  "
Set-TextCell 33 4 "1.75"
$ws.Cells.Item(33, 5).Value = "  -6.81%  "
Set-TextCell 34 4 "1.21"
$ws.Cells.Item(34, 5).Value = "  -8.60%  "
$ws.Cells.Item(35, 5).Value = "  +0.03%  "
Set-TextCell 36 4 "158.19"
$ws.Cells.Item(36, 5).Value = "  -1.38%  "
Set-TextCell 37 4 "0.113"
$ws.Cells.Item(37, 5).Value = "  -12.95%  "
Set-TextCell 38 4 "18.62"
$ws.Cells.Item(38, 5).Value = "  +0.23%  "
Set-TextCell 39 4 "18.15"
$ws.Cells.Item(39, 5).Value = "  -6.19%  "
Set-TextCell 40 4 "1.31"
$ws.Cells.Item(40, 5).Value = "  -9.56%  "
$ws.Cells.Item(41, 5).Value = "  +0.26%  "
$ws.Cells.Item(42, 5).Value = "  -7.28%  "
Set-TextCell 43 4 "4.76"
$ws.Cells.Item(43, 5).Value = "  -8.37%  "
Set-TextCell 44 4 "0.323"
$ws.Cells.Item(44, 5).Value = "  -8.21%  "
Set-TextCell 45 4 "2.34"
$ws.Cells.Item(45, 5).Value = "  -8.36%  "
Set-TextCell 46 4 "38.55"
$ws.Cells.Item(46, 5).Value = "  -2.32%  "
Set-TextCell 47 4 "140.17"
$ws.Cells.Item(47, 5).Value = "  -6.02%  "
Set-TextCell 48 4 "3.41"
$ws.Cells.Item(48, 5).Value = "  -8.68%  "
Set-TextCell 49 4 "0.506"
$ws.Cells.Item(49, 5).Value = "  -8.99%  "
$ws.Cells.Item(50, 4).Value = "0.0₆0251"
$ws.Cells.Item(50, 5).Value = "  -10.58%  "
Set-TextCell 51 4 "0.0728"
$ws.Cells.Item(51, 5).Value = "  -2.56%  "
